# Generate Report for Handoff
#
# Updates the localization-status workbook so each locale sheet (and the
# Overview roll-up) reflects that the handoff package is ready:
#   - Status cells move from "Handed back: in sync with en-US" to
#     "Ready for handoff"
#   - The associated timestamps are bumped forward a few seconds to the
#     new handoff-generation time
#   - The now-shorter status text no longer needs the wide status column,
#     so the Status columns are narrowed to match

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# --- Overview sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E2").Value = $newStatus
$ws.Range("F2").Value = $newStatus
$ws.Range("G2").Value = "2016-08-22 05:08:03"

# Status columns (zh-cn / de-de) on Overview were sized for the long
# "Handed back..." text; narrow them to fit "Ready for handoff".
$ws.Range("E1").ColumnWidth = 16.333333333333336
$ws.Range("F1").ColumnWidth = 16.333333333333336

# --- zh-cn sheet ------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C2").Value = $newStatus
$ws.Range("H2").Value = "2016-08-22 05:07:57"
$ws.Range("C1").ColumnWidth = 16.333333333333336

# --- de-de sheet ------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C2").Value = $newStatus
$ws.Range("H2").Value = "2016-08-22 05:08:03"
$ws.Range("C1").ColumnWidth = 16.333333333333336
